$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 764.5789
$ws.Range("I33").Value = 532.38464
$ws.Range("K33").Value = 532.38464
$ws.Range("M33").Value = -303.38464
$ws.Range("H98").Value = 4137.923
$ws.Range("I98").Value = 3072.182
$ws.Range("J98").Value = 9999.5
$ws.Range("K98").Value = 3072.182
$ws.Range("L98").Value = 9999.5
$ws.Range("M98").Value = -1574.182
$ws.Range("N98").Value = -12995.5
$ws.Range("H116").Value = 3203.1143
$ws.Range("I116").Value = 2468.3845
$ws.Range("J116").Value = 5325.6665
$ws.Range("K116").Value = 2468.3845
$ws.Range("L116").Value = 5325.6665
$ws.Range("M116").Value = 973.6154999999999
$ws.Range("N116").Value = -12209.6665
$ws.Range("H122").Value = 4137.923
$ws.Range("I122").Value = 3072.182
$ws.Range("J122").Value = 9999.5
$ws.Range("K122").Value = 9216.545999999998
$ws.Range("L122").Value = 29998.5
$ws.Range("M122").Value = -6766.545999999998
$ws.Range("N122").Value = -34898.5
$ws.Range("H129").Value = 4033690.2
$ws.Range("I129").Value = 31250300
$ws.Range("J129").Value = 1599.9073
$ws.Range("K129").Value = 93750900
$ws.Range("L129").Value = 4799.7219
$ws.Range("M129").Value = -93745900
$ws.Range("N129").Value = -14799.7219
$ws.Range("H133").Value = 28081.818
$ws.Range("J133").Value = 28081.818
$ws.Range("L133").Value = 28081.818
$ws.Range("N133").Value = -38201.818

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7029.6143
$ws.Range("I32").Value = 4322.5684
$ws.Range("J32").Value = 16191.923
$ws.Range("K32").Value = 4322.5684
$ws.Range("L32").Value = 16191.923
$ws.Range("M32").Value = -4035.5684
$ws.Range("N32").Value = -16765.923
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("H64").Value = 29850
$ws.Range("J64").Value = 29850
$ws.Range("L64").Value = 29850
$ws.Range("N64").Value = -30346
$ws.Range("H67").Value = 29850
$ws.Range("J67").Value = 29850
$ws.Range("L67").Value = 29850
$ws.Range("N67").Value = -31566
$ws.Range("H110").Value = 1108.5834
$ws.Range("I110").Value = 551.9032
$ws.Range("J110").Value = 4560
$ws.Range("K110").Value = 551.9032
$ws.Range("L110").Value = 4560
$ws.Range("M110").Value = 1493.0968
$ws.Range("N110").Value = -8650
$ws.Range("H122").Value = 2362.2144
$ws.Range("I122").Value = 1812
$ws.Range("J122").Value = 2774.875
$ws.Range("K122").Value = 5436
$ws.Range("L122").Value = 8324.625
$ws.Range("M122").Value = -2986
$ws.Range("N122").Value = -13224.625
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1675.3478
$ws.Range("I105").Value = 1426.25
$ws.Range("K105").Value = 1426.25
$ws.Range("M105").Value = 320.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2293.9473
$ws.Range("I122").Value = 2612.8572
$ws.Range("J122").Value = 2107.9167
$ws.Range("K122").Value = 7838.571599999999
$ws.Range("L122").Value = 6323.750100000001
$ws.Range("M122").Value = -5388.571599999999
$ws.Range("N122").Value = -11223.7501

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1270.5834
$ws.Range("I5").Value = 494.3125
$ws.Range("J5").Value = 2823.125
$ws.Range("K5").Value = 1482.9375
$ws.Range("L5").Value = 8469.375
$ws.Range("M5").Value = -1370.9375
$ws.Range("N5").Value = -8693.375
$ws.Range("H12").Value = 120.86207
$ws.Range("I12").Value = 16.5
$ws.Range("J12").Value = 148.08696
$ws.Range("K12").Value = 49.5
$ws.Range("L12").Value = 444.26088
$ws.Range("M12").Value = 123.5
$ws.Range("N12").Value = -790.26088
$ws.Range("H86").Value = 1166.6666
$ws.Range("H89").Value = 1166.6666
$ws.Range("H97").Value = 1600
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1600
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 4800
$ws.Range("N97").Value = -5792
$ws.Range("H107").Value = 903.92
$ws.Range("I107").Value = 1047.5
$ws.Range("J107").Value = 836.35297
$ws.Range("K107").Value = 3142.5
$ws.Range("L107").Value = 2509.05891
$ws.Range("M107").Value = -1222.5
$ws.Range("N107").Value = -6349.05891
$ws.Range("H131").Value = 1090.5577
$ws.Range("I131").Value = 973.4545000000001
$ws.Range("J131").Value = 1121.9756
$ws.Range("K131").Value = 2920.3635
$ws.Range("L131").Value = 3365.9268
$ws.Range("M131").Value = 2119.6365
$ws.Range("N131").Value = -13445.9268
$ws.Range("H135").Value = 1270.5834
$ws.Range("I135").Value = 494.3125
$ws.Range("J135").Value = 2823.125
$ws.Range("K135").Value = 4448.8125
$ws.Range("L135").Value = 25408.125
$ws.Range("M135").Value = -1913.8125
$ws.Range("N135").Value = -30478.125
$ws.Range("H138").Value = 910
$ws.Range("I138").Value = 910
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 2730
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 2410
$ws.Range("M97").ClearContents()
$ws.Range("N138").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("H122").Value = 4961.1113
$ws.Range("I122").Value = 11833.333
$ws.Range("J122").Value = 3586.6667
$ws.Range("K122").Value = 35499.999
$ws.Range("L122").Value = 10760.0001
$ws.Range("M122").Value = -33049.999
$ws.Range("N122").Value = -15660.0001
$ws.Range("H141").Value = 29000
$ws.Range("J141").Value = 29000
$ws.Range("L141").Value = 29000
$ws.Range("N141").Value = -39360
$ws.Range("N101").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 698.5714
$ws.Range("I16").Value = 731.6667
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 731.6667
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -561.6667
$ws.Range("N16").Value = -840
$ws.Range("H61").Value = 100002030
$ws.Range("I61").Value = 142858340
$ws.Range("J61").Value = 3998.3333
$ws.Range("K61").Value = 142858340
$ws.Range("L61").Value = 3998.3333
$ws.Range("M61").Value = -142858138
$ws.Range("N61").Value = -4402.3333
$ws.Range("H68").Value = 1470.1305
$ws.Range("I68").Value = 1020
$ws.Range("J68").Value = 3090.6
$ws.Range("K68").Value = 1020
$ws.Range("L68").Value = 3090.6
$ws.Range("M68").Value = -271
$ws.Range("N68").Value = -4588.6
$ws.Range("H71").Value = 1470.1305
$ws.Range("I71").Value = 1020
$ws.Range("J71").Value = 3090.6
$ws.Range("K71").Value = 5100
$ws.Range("L71").Value = 15453
$ws.Range("M71").Value = -1356
$ws.Range("N71").Value = -22941
$ws.Range("H113").Value = 100002030
$ws.Range("I113").Value = 142858340
$ws.Range("J113").Value = 3998.3333
$ws.Range("K113").Value = 142858340
$ws.Range("L113").Value = 3998.3333
$ws.Range("M113").Value = -142856170
$ws.Range("N113").Value = -8338.3333
$ws.Range("H122").Value = 2628.9092
$ws.Range("I122").Value = 2437.6428
$ws.Range("K122").Value = 7312.928400000001
$ws.Range("M122").Value = -4862.928400000001
$ws.Range("H123").Value = 20000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 20000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -29800
$ws.Range("H136").Value = 5885252.5
$ws.Range("I136").Value = 20001880
$ws.Range("J136").Value = 3324.1667
$ws.Range("K136").Value = 60005640
$ws.Range("L136").Value = 9972.500100000001
$ws.Range("M136").Value = -60003090
$ws.Range("N136").Value = -15072.5001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 28338.334
$ws.Range("J75").Value = 28338.334
$ws.Range("L75").Value = 28338.334
$ws.Range("N75").Value = -30210.334
$ws.Range("H78").Value = 28338.334
$ws.Range("J78").Value = 28338.334
$ws.Range("L78").Value = 85015.00199999999
$ws.Range("N78").Value = -94375.00199999999
$ws.Range("H81").Value = 1035
$ws.Range("I81").Value = 745.25
$ws.Range("J81").Value = 1324.75
$ws.Range("K81").Value = 1490.5
$ws.Range("L81").Value = 2649.5
$ws.Range("M81").Value = -429.5
$ws.Range("N81").Value = -4771.5
$ws.Range("H84").Value = 1035
$ws.Range("I84").Value = 745.25
$ws.Range("J84").Value = 1324.75
$ws.Range("K84").Value = 7452.5
$ws.Range("L84").Value = 13247.5
$ws.Range("M84").Value = -2148.5
$ws.Range("N84").Value = -23855.5
$ws.Range("H106").Value = 29000
$ws.Range("J106").Value = 29000
$ws.Range("L106").Value = 29000
$ws.Range("N106").Value = -31524
$ws.Range("H109").Value = 36500
$ws.Range("J109").Value = 36500
$ws.Range("L109").Value = 36500
$ws.Range("N109").Value = -39274
$ws.Range("H126").Value = 5883971.5
$ws.Range("I126").Value = 1026
$ws.Range("J126").Value = 11113256
$ws.Range("K126").Value = 3078
$ws.Range("L126").Value = 33339768
$ws.Range("M126").Value = -608
$ws.Range("N126").Value = -33344708
$ws.Range("H132").Value = 304305.8
$ws.Range("I132").Value = 481366.66
$ws.Range("J132").Value = 38714.5
$ws.Range("K132").Value = 1444099.98
$ws.Range("L132").Value = 116143.5
$ws.Range("M132").Value = -1441569.98
$ws.Range("N132").Value = -121203.5
$ws.Range("H136").Value = 2547.077
$ws.Range("I136").Value = 1630.2858
$ws.Range("J136").Value = 3616.6667
$ws.Range("K136").Value = 4890.857400000001
$ws.Range("L136").Value = 10850.0001
$ws.Range("M136").Value = -2340.857400000001
$ws.Range("N136").Value = -15950.0001
$ws.Range("H138").Value = 29833.334
$ws.Range("J138").Value = 29833.334
$ws.Range("L138").Value = 29833.334
$ws.Range("N138").Value = -40113.334
